$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# --- Row 10 ("No.") ---
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("A10").HorizontalAlignment = $xlCenter
$ws.Range("B10").Value = 9
$ws.Range("D10").Value = 19
$ws.Range("E10").Value = 28

# --- Row 11 ("Marking") ---
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("A11").HorizontalAlignment = $xlCenter
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# --- Row 12 ("Total") ---
$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("A12").HorizontalAlignment = $xlCenter
$ws.Range("B12").Value = 36
$ws.Range("E12").Value = "36/112"

# --- Row 15: drop the third (Student Ans / Correct Ans) header pair ---
$ws.Range("G15:H15").Clear()

# --- Mark some "Student Ans" (column A) cells as matching the correct
#     answer, using the green "correctStyle" used elsewhere in the sheet ---
$ws.Range("A16").Style = "correctStyle"
$ws.Range("A16").HorizontalAlignment = $xlCenter
$ws.Range("A16").Value = "Option A"

$ws.Range("A18").Style = "correctStyle"
$ws.Range("A18").HorizontalAlignment = $xlCenter
$ws.Range("A18").Value = "Option B"

$ws.Range("A19").Style = "correctStyle"
$ws.Range("A19").HorizontalAlignment = $xlCenter
$ws.Range("A19").Value = "Option C"

$ws.Range("A25").Style = "correctStyle"
$ws.Range("A25").HorizontalAlignment = $xlCenter
$ws.Range("A25").Value = "Option A"

$ws.Range("A27").Style = "correctStyle"
$ws.Range("A27").HorizontalAlignment = $xlCenter
$ws.Range("A27").Value = "Option A"

$ws.Range("A32").Style = "correctStyle"
$ws.Range("A32").HorizontalAlignment = $xlCenter
$ws.Range("A32").Value = "Option C"

$ws.Range("A33").Style = "correctStyle"
$ws.Range("A33").HorizontalAlignment = $xlCenter
$ws.Range("A33").Value = "Option D"

# --- Column D of the second (Student Ans / Correct Ans) pair only keeps
#     data in rows 16 and 18; everything else below is cleared ---
$ws.Range("D16").Style = "correctStyle"
$ws.Range("D16").HorizontalAlignment = $xlCenter
$ws.Range("D16").Value = "Option A"

$ws.Range("D18").Style = "correctStyle"
$ws.Range("D18").HorizontalAlignment = $xlCenter
$ws.Range("D18").Value = "Option D"

# --- Remove the now-unused D/E data (rows 19-40) and all of the third
#     Student Ans / Correct Ans column pair (G/H, rows 16-40) ---
$ws.Range("D19:E40").Clear()
$ws.Range("G16:H40").Clear()
